$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 11, shifting existing rows 11-54 down to 12-55.
$ws.Rows.Item(11).Insert()

# Populate the newly inserted row with the new label data.
$ws.Range("A11").Value = "a1.9"
$ws.Range("B11").Value = "b1.9"
$ws.Range("C11").Value = "Personnel"
$ws.Range("D11").Value = "Forwarding to the right area while getting row ticket"

# Reflect the active selection that was left on D11 after the edit.
$ws.Range("D11").Select() | Out-Null
